$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 416
$firstRow = 2

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45206
